# "reorganize Ar4 and Ar3" - rebuild Sheet1's "head" table into the new
# "right_arm" table (Ar4 arm diagram data): new headers/column order, new
# row data (right shoulder/bicept joints), and drop the two trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Headers (row 1): Name | Num | Degree (A) | Min | Mid | Max ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Num"
$ws.Range("C1").Value = "Degree (A)"
$ws.Range("D1").Value = "Min"
$ws.Range("E1").Value = "Mid"
$ws.Range("F1").Value = "Max"

# --- Data rows 2-5: Name | Num | Degree(A)=0..3, Min/Mid/Max left blank ---
$ws.Range("A2").Value = "right_shoulder_x"
$ws.Range("B2").Value = "M4"
$ws.Range("C2").Value = 0
$ws.Range("D2:F2").ClearContents()

$ws.Range("A3").Value = "right_shoulder_y"
$ws.Range("B3").Value = "M3"
$ws.Range("C3").Value = 1
$ws.Range("D3:F3").ClearContents()

$ws.Range("A4").Value = "right_shoulder_z"
$ws.Range("B4").Value = "M1"
$ws.Range("C4").Value = 2
$ws.Range("D4:F4").ClearContents()

$ws.Range("A5").Value = "right_bicept"
$ws.Range("B5").Value = "M2"
$ws.Range("C5").Value = 3
$ws.Range("D5:F5").ClearContents()

# Rows 6 and 7 no longer exist in the reorganized table - remove them.
$ws.Rows("6:7").Delete()

# Shrink the table (and its autofilter) from A1:F7 down to A1:F5, and
# rename it from "head" to "right_arm".
$lo.Resize($ws.Range("A1:F5"))
$lo.Name = "right_arm"

# Column C ("Degree (A)") gets a custom width in the new layout.
$ws.Columns("C").ColumnWidth = 11.5

# Selection moves to G15 in the saved view.
$ws.Range("G15").Select()
